# Append the two "restart" log rows (car reset coordinate + animation stop
# events) that were missing from the results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Tue_Dec__5_09_17_31_2023"
$ws.Range("B3").Value = "f"
$ws.Range("C3").Value = 30

$ws.Range("A4").Value = "Tue_Dec__5_09_20_16_2023"
$ws.Range("B4").Value = "f"
$ws.Range("C4").Value = 30
